$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '65.402.45'
$ws.Range("E2").Value = '  -2.17%  '

# Row 3
Set-TextValue $ws.Range("D3") '3.384.41'
$ws.Range("E3").Value = '  -2.53%  '

# Row 4
Set-TextValue $ws.Range("D4") '1.00'
$ws.Range("E4").Value = '  +0.07%  '

# Row 5
Set-TextValue $ws.Range("D5") '594.54'
$ws.Range("E5").Value = '  -1.69%  '

# Row 6
Set-TextValue $ws.Range("D6") '141.00'
$ws.Range("E6").Value = '  -5.13%  '

# Row 7
$ws.Range("E7").Value = '  -0.06%  '

# Row 8
Set-TextValue $ws.Range("D8") '3.381.12'
$ws.Range("E8").Value = '  -2.57%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.468'
$ws.Range("E9").Value = '  -3.21%  '

# Row 10
Set-TextValue $ws.Range("D10") '7.92'
$ws.Range("E10").Value = '  +4.57%  '

# Row 11
$ws.Range("E11").Value = '  -6.97%  '

# Row 12
Set-TextValue $ws.Range("D12") '0.405'
$ws.Range("E12").Value = '  -4.79%  '

# Row 13
Set-TextValue $ws.Range("D13") '3.957.59'
$ws.Range("E13").Value = '  -2.51%  '

# Row 14
$ws.Range("E14").Value = '  -7.56%  '

# Row 15
Set-TextValue $ws.Range("D15") '29.49'
$ws.Range("E15").Value = '  -7.39%  '

# Row 16
$ws.Range("E16").Value = '  -0.66%  '

# Row 17
Set-TextValue $ws.Range("D17") '65.375.83'
$ws.Range("E17").Value = '  -2.16%  '

# Row 18
Set-TextValue $ws.Range("D18") '3.382.81'
$ws.Range("E18").Value = '  -2.61%  '

# Row 19
Set-TextValue $ws.Range("D19") '10.40'
$ws.Range("E19").Value = '  +2.73%  '

# Row 20
Set-TextValue $ws.Range("D20") '6.09'
$ws.Range("E20").Value = '  -6.10%  '

# Row 21
Set-TextValue $ws.Range("D21") '14.57'
$ws.Range("E21").Value = '  -5.67%  '

# Row 22
Set-TextValue $ws.Range("D22") '413.06'
$ws.Range("E22").Value = '  -6.27%  '

# Row 23
Set-TextValue $ws.Range("D23") '0.577'
$ws.Range("E23").Value = '  -5.82%  '

# Row 24
Set-TextValue $ws.Range("D24") '77.08'
$ws.Range("E24").Value = '  -2.88%  '

# Row 25
$ws.Range("E25").Value = '  +0.04%  '

# Row 26
Set-TextValue $ws.Range("D26") '3.518.95'
$ws.Range("E26").Value = '  -2.48%  '

# Row 27
$ws.Range("E27").Value = '  -10.18%  '

# Row 28
Set-TextValue $ws.Range("D28") '9.19'
$ws.Range("E28").Value = '  -6.30%  '

# Row 29
Set-TextValue $ws.Range("D29") '7.73'
$ws.Range("E29").Value = '  -7.92%  '

# Row 30
$ws.Range("E30").Value = '  -3.26%  '

# Row 31
$ws.Range("E31").Value = '  +0.26%  '

# Row 32
$ws.Range("E32").Value = '  -5.24%  '

# Row 33
$ws.Range("E33").Value = '  -8.83%  '

# Row 34
Set-TextValue $ws.Range("D34") '24.29'
$ws.Range("E34").Value = '  -4.60%  '

# Row 35
Set-TextValue $ws.Range("D35") '3.381.47'
$ws.Range("E35").Value = '  -2.32%  '

# Row 37
Set-TextValue $ws.Range("D37") '5.51'
$ws.Range("E37").Value = '  -9.40%  '

# Row 38
Set-TextValue $ws.Range("D38") '1.67'
$ws.Range("E38").Value = '  -7.40%  '

# Row 39
Set-TextValue $ws.Range("D39") '7.48'
$ws.Range("E39").Value = '  -5.82%  '

# Row 40
Set-TextValue $ws.Range("D40") '1.00'
$ws.Range("E40").Value = '  +0.08%  '

# Row 41
Set-TextValue $ws.Range("D41") '167.22'
$ws.Range("E41").Value = '  -5.57%  '

# Row 42
Set-TextValue $ws.Range("D42") '0.0851'
$ws.Range("E42").Value = '  -4.51%  '

# Row 43
Set-TextValue $ws.Range("D43") '0.868'
$ws.Range("E43").Value = '  -2.20%  '

# Row 44
Set-TextValue $ws.Range("D44") '5.00'
$ws.Range("E44").Value = '  -7.96%  '

# Row 45
$ws.Range("E45").Value = '  -11.27%  '

# Row 46
$ws.Range("E46").Value = '  -2.04%  '

# Row 47
Set-TextValue $ws.Range("D47") '26.54'
$ws.Range("E47").Value = '  -9.36%  '

# Row 48
$ws.Range("E48").Value = '  -5.60%  '

# Row 49
$ws.Range("E49").Value = '  -6.18%  '

# Row 50
Set-TextValue $ws.Range("D50") '2.24'
$ws.Range("E50").Value = '  -9.12%  '

# Row 51
Set-TextValue $ws.Range("B51") 'SuiNetwork'
Set-TextValue $ws.Range("C51") 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue $ws.Range("D51") '0.911'
$ws.Range("E51").Value = '  -7.91%  '
